# Apply updated crypto price/volume figures (Price column D, Volume(1h) column E).
# Values are stored as plain text in the sheet, so numeric-looking Price
# entries are written with a leading apostrophe to keep Excel from
# reinterpreting them as numbers (which would change their displayed form).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.596.99"
$ws.Range("E2").Value = "  +3.60%  "
$ws.Range("D3").Value = "1.919.21"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'248.79"
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("D6").Value = "'0.695"
$ws.Range("E6").Value = "  +1.30%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'44.19"
$ws.Range("E8").Value = "  +1.91%  "
$ws.Range("D9").Value = "'58.53"
$ws.Range("E9").Value = "  +9.60%  "
$ws.Range("E10").Value = "  +3.61%  "
$ws.Range("E11").Value = "  +3.52%  "
$ws.Range("D12").Value = "'0.0996"
$ws.Range("E12").Value = "  +2.63%  "
$ws.Range("D13").Value = "'14.56"
$ws.Range("E13").Value = "  +8.35%  "
$ws.Range("D14").Value = "'0.797"
$ws.Range("E14").Value = "  +4.76%  "
$ws.Range("D15").Value = "2.201.72"
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("E16").Value = "  +4.91%  "
$ws.Range("D17").Value = "1.922.82"
$ws.Range("E17").Value = "  +2.23%  "
$ws.Range("D18").Value = "36.581.01"
$ws.Range("E18").Value = "  +3.17%  "
$ws.Range("D19").Value = "'74.27"
$ws.Range("E19").Value = "  +1.96%  "
$ws.Range("D20").Value = "0.0₃0859"
$ws.Range("E20").Value = "  +4.80%  "
$ws.Range("E21").Value = "  +3.04%  "
$ws.Range("D22").Value = "'13.23"
$ws.Range("E22").Value = "  +3.72%  "
$ws.Range("D23").Value = "'5.20"
$ws.Range("E23").Value = "  +5.26%  "
$ws.Range("E24").Value = "  +2.13%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").Value = "'2.19"
$ws.Range("E26").Value = "  +1.51%  "
$ws.Range("D27").Value = "'167.63"
$ws.Range("E27").Value = "  +1.55%  "
$ws.Range("D28").Value = "'8.82"
$ws.Range("E28").Value = "  +3.42%  "
$ws.Range("D29").Value = "'18.76"
$ws.Range("E29").Value = "  +2.45%  "
$ws.Range("E30").Value = "  +1.74%  "
$ws.Range("D31").Value = "'4.54"
$ws.Range("D32").Value = "'0.0609"
$ws.Range("E32").Value = "  +3.88%  "
$ws.Range("D33").Value = "'1.99"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  +5.04%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "'0.0852"
$ws.Range("E36").Value = "  +18.70%  "
$ws.Range("D37").Value = "'1.50"
$ws.Range("E37").Value = "  -13.58%  "
$ws.Range("D38").Value = "'0.881"
$ws.Range("E38").Value = "  +4.42%  "
$ws.Range("D39").Value = "'17.69"
$ws.Range("E39").Value = "  +45.50%  "
$ws.Range("D40").Value = "'2.01"
$ws.Range("E40").Value = "  +3.42%  "
$ws.Range("D41").Value = "'107.04"
$ws.Range("E41").Value = "  +11.78%  "
$ws.Range("D42").Value = "'0.0228"
$ws.Range("E42").Value = "  +5.11%  "
$ws.Range("D43").Value = "'17.16"
$ws.Range("E43").Value = "  -1.47%  "
$ws.Range("E44").Value = "  +3.11%  "
$ws.Range("D45").Value = "1.339.23"
$ws.Range("E45").Value = "  +2.75%  "
$ws.Range("D46").Value = "'2.36"
$ws.Range("E46").Value = "  +0.92%  "
$ws.Range("D47").Value = "'2.49"
$ws.Range("E47").Value = "  +4.51%  "
$ws.Range("D48").Value = "'0.0812"
$ws.Range("E48").Value = "  +1.96%  "
$ws.Range("D49").Value = "'2.80"
$ws.Range("E49").Value = "  +2.78%  "
$ws.Range("D50").Value = "'6.43"
$ws.Range("E50").Value = "  +3.25%  "
$ws.Range("D51").Value = "2.107.24"
$ws.Range("E51").Value = "  +2.12%  "
